$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.278.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -7.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.534.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.96%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "391.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "122.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.66%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.524.87"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.98%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.582"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -10.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.675"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -12.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -22.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000323"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -25.68%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.085.66"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.53%  "

$ws.Range("E16").Value = "  -3.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.523.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "63.229.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "390.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -14.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.06"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.24%  "

$ws.Range("E26").Value = "  +10.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.80%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -14.64%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.47%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.07%  "

$ws.Range("E32").Value = "  -8.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.16%  "

$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.148"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.56%  "

$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "36.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0433"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -11.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.997"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0644"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -18.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.130"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -13.40%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +20.49%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "139.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.91%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +15.09%  "

$ws.Range("E46").Value = "  -0.22%  "

$ws.Range("B47").Value = "LidoDAOToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.15%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.85%  "

$ws.Range("E49").Value = "  -9.79%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -10.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.272"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.68%  "
